$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whitespace from the street-name values in column N (rows 2-5)
$ws.Range("N2").Value = "RoermonderStr.79"
$ws.Range("N3").Value = "Bachstraße4"
$ws.Range("N4").Value = "Muehlenstraße21f"
$ws.Range("N5").Value = "VaalserStr.226"

# Update the view state: scroll so column C is the left-most visible column,
# and move the active selection to M7
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
$ws.Range("M7").Select()
